# Re-create the "_GoBack" bookmark that Word silently (re-)inserts at the
# very start of the document content whenever it is opened and saved again
# (it marks the last editing position for the Ctrl+Alt+Z "go back" feature).
#
# The target location is immediately after the first paragraph's <w:pPr>
# and before its first run, i.e. a zero-length bookmark that wraps the
# very first character position of the document.
#
# A bookmark added directly at the absolute document start (position 0)
# behaves oddly (its start/end markers get attached to different runs
# across the first paragraph boundary), so we nudge it one character in:
# insert a one-character placeholder at the start, anchor the empty
# bookmark right after it (still effectively "at the start" once the
# placeholder is removed), then delete the placeholder again.

$d = $word.ActiveDocument

$docStart = $d.Content.Start

# Insert a throwaway placeholder character at the very beginning.
$headRange = $d.Range($docStart, $docStart)
$headRange.InsertBefore("X")

# Anchor the bookmark right after the placeholder (collapsed / zero length).
$markRange = $d.Range($docStart + 1, $docStart + 1)
$d.Bookmarks.Add("_GoBack", $markRange)

# Remove the placeholder character again, leaving the bookmark at the
# true start of the document.
$d.Range($docStart, $docStart + 1).Delete()
